# Extend the grid from columns A:D to A:F (two more "video" columns),
# mirroring the pattern already used for columns C/D:
#   - header row 1 gets sequential numbers continuing 0,1,2 -> 3,4 in E1/F1
#   - any data row that already has a C value gets the same value copied to D
#     (if D is still empty) and, for rows that already span out to D (C+D),
#     the same pair is copied again into E/F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: copy the look of D1 into the two new header cells, then set values ---
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4

# --- Rows whose only "data" column is C: duplicate the value into D ---
$rowsCtoD = @(3, 10, 12, 14, 22, 24, 27, 51, 76, 96, 106, 161)
foreach ($r in $rowsCtoD) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 4).Value = $cVal
}

# --- Rows that already have both C and D filled: duplicate the C/D pair into E/F ---
$rowsCDtoEF = @(119, 145)
foreach ($r in $rowsCDtoEF) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 5).Value = $cVal
    $ws.Cells.Item($r, 6).Value = $dVal
}

Write-Host "Updated dimension: $($ws.UsedRange.Address())"
